$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM_PERMISSION")

$ws.Range("G8").Value = "/notify/send-setting"
$ws.Range("G9").Value = "/notify/mail-template"
$ws.Range("G10").Value = "/notify/inmail-template"
$ws.Range("G11").Value = "/notify/mail-setting"
$ws.Range("G12").Value = "/notify/msg-record"
$ws.Range("G13").Value = "/notify/announcement"
$ws.Range("G14").Value = "/notify/user-msg"
$ws.Range("G15").Value = "/notify/receive-setting"
